# "Updated view manage order script"
#
# The ManageProducts QA fixture keeps a handful of randomly generated
# "prodXXXX" product-name placeholders in column B of the Input sheet.
# This refreshes the three that the "view manage order" test script reads
# (rows 2, 3 and 5) with freshly generated product names, matching the
# regenerated fixture data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "prodAqci"
$ws.Range("B3").Value = "prodrpbv"
$ws.Range("B5").Value = "prodqGOY"
